# Weekly update: insert the latest week's record as a new row 66 (ahead of
# the existing history), pushing the older rows down by one. This matches
# the "logica_diaria" -> weekly snapshot pattern used by this subconjunto.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 66; everything from 66..100
# shifts down to 67..101 (dimension grows to A1:T101 automatically).
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with this week's record.
$ws.Cells.Item(66, 1).Value  = 5
$ws.Cells.Item(66, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(66, 3).Value  = "Maule"
$ws.Cells.Item(66, 4).Value  = 44572
$ws.Cells.Item(66, 5).Value  = 7
$ws.Cells.Item(66, 6).Value  = "Fruta"
$ws.Cells.Item(66, 7).Value  = 100108
$ws.Cells.Item(66, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(66, 9).Value  = 100108002
$ws.Cells.Item(66, 10).Value = "Mango"
$ws.Cells.Item(66, 11).Value = "Sin especificar"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 170
$ws.Cells.Item(66, 14).Value = 6000
$ws.Cells.Item(66, 15).Value = 6000
$ws.Cells.Item(66, 16).Value = 6000
$ws.Cells.Item(66, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(66, 18).Value = "Perú"
$ws.Cells.Item(66, 19).Value = 1500
$ws.Cells.Item(66, 20).Value = 4
